# chore: update Sheets via scheduled runner
#
# Applies updated market-price / profit figures to several leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 19997.5
$ws.Range("J13").Value = 19997.5
$ws.Range("L13").Value = 19997.5
$ws.Range("N13").Value = -20335.5

$ws.Range("H64").Value = 3680
$ws.Range("I64").Value = 3680
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3680
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3432
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 3680
$ws.Range("I67").Value = 3680
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3680
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2822
$ws.Range("N67").ClearContents()

# --- ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 8700
$ws.Range("J9").Value = 8700
$ws.Range("L9").Value = 8700
$ws.Range("N9").Value = -9040

$ws.Range("H20").Value = 8700
$ws.Range("J20").Value = 8700
$ws.Range("L20").Value = 8700
$ws.Range("N20").Value = -9240

$ws.Range("H45").Value = 250000740
$ws.Range("I45").Value = 250000740
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 250000740
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -250000363
$ws.Range("N45").ClearContents()

$ws.Range("H63").Value = 2673.2812
$ws.Range("I63").Value = 2716.5557
$ws.Range("J63").Value = 2439.6
$ws.Range("K63").Value = 2716.5557
$ws.Range("L63").Value = 2439.6
$ws.Range("M63").Value = -2030.5557
$ws.Range("N63").Value = -3811.6

$ws.Range("H66").Value = 2673.2812
$ws.Range("I66").Value = 2716.5557
$ws.Range("J66").Value = 2439.6
$ws.Range("K66").Value = 13582.7785
$ws.Range("L66").Value = 12198
$ws.Range("M66").Value = -10150.7785
$ws.Range("N66").Value = -19062

$ws.Range("H88").Value = 2324.0476
$ws.Range("I88").Value = 1329.2858
$ws.Range("J88").Value = 2821.4285
$ws.Range("K88").Value = 1329.2858
$ws.Range("L88").Value = 2821.4285
$ws.Range("M88").Value = -923.2858000000001
$ws.Range("N88").Value = -3633.4285

$ws.Range("H91").Value = 2324.0476
$ws.Range("I91").Value = 1329.2858
$ws.Range("J91").Value = 2821.4285
$ws.Range("K91").Value = 1329.2858
$ws.Range("L91").Value = 2821.4285
$ws.Range("M91").Value = 74.71419999999989
$ws.Range("N91").Value = -5629.4285

$ws.Range("H94").Value = 27040
$ws.Range("J94").Value = 27040
$ws.Range("L94").Value = 27040
$ws.Range("N94").Value = -28842

$ws.Range("H101").Value = 52864
$ws.Range("J101").Value = 52864
$ws.Range("L101").Value = 52864
$ws.Range("N101").Value = -59354

# --- BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 40000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 40000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 40000
$ws.Range("N15").Value = -40454
$ws.Range("M15").ClearContents()

$ws.Range("H19").Value = 26900
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 26900
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 26900
$ws.Range("N19").Value = -27246
$ws.Range("M19").ClearContents()

$ws.Range("H92").Value = 28439.8
$ws.Range("J92").Value = 28439.8
$ws.Range("L92").Value = 28439.8
$ws.Range("N92").Value = -33431.8

$ws.Range("H107").Value = 1217.4375
$ws.Range("I107").Value = 1040.4286
$ws.Range("K107").Value = 1040.4286
$ws.Range("M107").Value = 879.5714

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 3880
$ws.Range("I25").Value = 3880
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 3880
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -3706
$ws.Range("N25").ClearContents()

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1076.3636
$ws.Range("I108").Value = 480
$ws.Range("J108").Value = 2666.6667
$ws.Range("K108").Value = 1440
$ws.Range("L108").Value = 8000.000100000001
$ws.Range("M108").Value = 1440
$ws.Range("N108").Value = -13760.0001

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 8114.2856
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 8114.2856
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 8114.2856
$ws.Range("N18").Value = -8700.285599999999
$ws.Range("M18").ClearContents()

$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492

$ws.Range("H105").Value = 13672
$ws.Range("J105").Value = 13672
$ws.Range("L105").Value = 13672
$ws.Range("N105").Value = -20660

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 9011.111000000001
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 9887.5
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 9887.5
$ws.Range("M23").Value = -1770
$ws.Range("N23").Value = -10347.5

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H30").Value = 33333.332
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 33333.332
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 33333.332
$ws.Range("N30").Value = -33549.332
$ws.Range("M30").ClearContents()

$ws.Range("H55").Value = 239.78947
$ws.Range("I55").Value = 175.875
$ws.Range("J55").Value = 580.6667
$ws.Range("K55").Value = 175.875
$ws.Range("L55").Value = 580.6667
$ws.Range("M55").Value = -2.875
$ws.Range("N55").Value = -926.6667

$ws.Range("H68").Value = 9957.143
$ws.Range("I68").Value = 100000
$ws.Range("J68").Value = 3030.7693
$ws.Range("K68").Value = 100000
$ws.Range("L68").Value = 3030.7693
$ws.Range("M68").Value = -99251
$ws.Range("N68").Value = -4528.7693

$ws.Range("H71").Value = 9957.143
$ws.Range("I71").Value = 100000
$ws.Range("J71").Value = 3030.7693
$ws.Range("K71").Value = 500000
$ws.Range("L71").Value = 15153.8465
$ws.Range("M71").Value = -496256
